$d = $word.ActiveDocument

# 1. Remove the _GoBack bookmark from its current position (around "{CI}" in
#    the "C.I.:" paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Make the "{telefonos_de_contacto}" paragraph bold (applies to both the
#    paragraph-mark's run properties and every run's run properties).
$pTel = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*telefonos_de_contacto*") {
        $pTel = $cand
        break
    }
}
$pTel.Range.Font.Bold = 1

# 3. Re-insert the (now empty-range) _GoBack bookmark inside the word
#    "FIRMADO" of the consent paragraph, splitting the run there into
#    "...SER F" and "IRMADO POR UNO DE SUS PADRES)".
$rngFind = $d.Content
$rngFind.Find.ClearFormatting()
$found = $rngFind.Find.Execute("IRMADO", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPos = $rngFind.Start

$pLast = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*FIRMADO POR UNO*") {
        $pLast = $cand
        break
    }
}
$tail = $d.Range($splitPos, $pLast.Range.End - 1)
$tail.Font.Bold = 0
$tail.Font.Bold = 1

$d.Bookmarks.Add("_GoBack", $d.Range($splitPos, $splitPos))
